# Update "想去人数" (attendee interest count) values in the "展览" and
# "全部类型" worksheets to reflect freshly scraped data.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 201
$ws1.Range("F4").Value = 804
$ws1.Range("F6").Value = 19

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 201
$ws4.Range("F5").Value = 804
$ws4.Range("F7").Value = 19
